$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing observed-data rows (2-14) down by one row to make room
# for a new "Day 0" row, writing from the bottom up so later rows keep their
# final values. Row 3's value is also updated (956949.1003628989 ->
# 456949.10036289901) per the new Day-0/Day-1 inputs.

$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 2).Value = 6403050.0364428777

$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 2).Value = 6475788.1657214472

$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = 5693093.0453154361

$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = 5355557.2857550187

$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = 4512520.3528595893

$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = 4113185.3448669091

$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = 3166443.706749748

$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = 2984699.4576945822

$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(7, 2).Value = 2828412.8733981312

$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = 2993657.5240874859

$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = 1723601.156960699

$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = 820110.43364600511

$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = 456949.10036289901

$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = 300000

# Update the active selection, as left by the user after the edit.
$ws.Range("M27").Select()
